$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 is an existing data row; the new row 3 duplicates it exactly.
# Copy row 2 (A2:L2) and paste only the values (not formats/styles) into row 3
# so the new cells keep the default (unstyled) look, matching the target.
$ws.Range("A2:L2").Copy() | Out-Null
$ws.Range("A3").PasteSpecial(-4163) | Out-Null
$excel.CutCopyMode = 0
